$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.826271637843651
$ws.Range("C3").Value = 1.906592789843651
$ws.Range("C4").Value = 1.910950668843651
$ws.Range("C5").Value = 2.452713729843651
$ws.Range("C6").Value = 1.806423198843651
$ws.Range("B7").Value = 1.454170430843651
$ws.Range("B8").Value = 1.527312638843651
$ws.Range("B9").Value = -0.2024072381563491
$ws.Range("D46").Value = 0.4006729710360878
$ws.Range("D47").Value = 0.3913241300360878
$ws.Range("D48").Value = 0.07541333303608777
$ws.Range("D49").Value = 0.1149542130360878
$ws.Range("C50").Value = 0.7012562470360878
$ws.Range("C51").Value = 1.029332457036088
$ws.Range("C52").Value = 1.485428559036088
$ws.Range("C53").Value = 1.831943276036088
$ws.Range("D53").Value = 0.55793444658209
$ws.Range("C54").Value = 1.696990870036088
$ws.Range("D54").Value = 0.44417418258209
$ws.Range("C55").Value = 0.1642582080360878
$ws.Range("D55").Value = -0.54231891241791
$ws.Range("B56").Value = 0.1757658360360879
$ws.Range("C56").Value = -0.63179670641791
$ws.Range("B57").Value = 0.1903804690360879
$ws.Range("C57").Value = 0.24128467758209
$ws.Range("C58").Value = 0.44433217858209
$ws.Range("C59").Value = 0.6062157845820899
$ws.Range("D59").Value = -0.8596988570317647
$ws.Range("C60").Value = 1.06983264558209
$ws.Range("D60").Value = -0.9590976180317647
$ws.Range("C61").Value = 0.7401392785820899
$ws.Range("D61").Value = -1.043273225031765
$ws.Range("B62").Value = 0.65042024058209
$ws.Range("C62").Value = -0.5948165710317647
$ws.Range("B63").Value = 0.66923223158209
$ws.Range("C63").Value = -1.459370673031765
$ws.Range("C64").Value = -1.543924081031765
$ws.Range("C65").Value = -1.197149625031765
$ws.Range("C66").Value = -0.8452485510317647
$ws.Range("C67").Value = -1.137602302031765
$ws.Range("B68").Value = -1.037127551031765
$ws.Range("B69").Value = -1.085847230031765
$ws.Range("D71").Value = 0.7977233685636995
$ws.Range("D72").Value = 1.031128409563699
$ws.Range("D73").Value = 0.9920349685636995
$ws.Range("D74").Value = 1.0168253855637
$ws.Range("C75").Value = 1.053976176563699
$ws.Range("C76").Value = 1.093632598563699
$ws.Range("C77").Value = 1.124545940563699
$ws.Range("C78").Value = 1.0707846905637
$ws.Range("D78").Value = 0.8985202055291455
$ws.Range("C79").Value = 0.9333254135636995
$ws.Range("D79").Value = 0.8577491755291454
$ws.Range("C80").Value = 0.9579434665636994
$ws.Range("D80").Value = 0.8987580675291454
$ws.Range("B81").Value = 1.0751219075637
$ws.Range("C81").Value = 1.258920946529146
$ws.Range("B82").Value = 1.1210615525637
$ws.Range("C82").Value = 1.493423065529146
$ws.Range("C83").Value = 1.027544699529146
$ws.Range("C84").Value = 0.5469191505291455
$ws.Range("D84").Value = -0.9817231827224345
$ws.Range("C85").Value = 0.6114425455291455
$ws.Range("D85").Value = -0.9854454907224345
$ws.Range("C86").Value = 0.3885492905291454
$ws.Range("D86").Value = -0.8873425837224346
$ws.Range("C87").Value = -0.1343109554708546
$ws.Range("D87").Value = -1.085683986722434
$ws.Range("B88").Value = 0.4230266025291455
$ws.Range("C88").Value = 0.1374695422775655
$ws.Range("B89").Value = 0.4182579295291455
$ws.Range("C89").Value = -0.2331624027224345
$ws.Range("C90").Value = -0.03326235772243452
$ws.Range("C91").Value = -0.2778009377224345
$ws.Range("D91").Value = -1.165608299123972
$ws.Range("C92").Value = -0.6454920347224345
$ws.Range("D92").Value = -1.024285167559779
$ws.Range("C93").Value = -0.7103774527224346
$ws.Range("D93").Value = -1.000794259827642
$ws.Range("B94").Value = -0.5011116027224345
$ws.Range("C94").Value = -0.6420067559859775
$ws.Range("B95").Value = -0.5544081717224345
$ws.Range("C95").Value = -0.7670271480143395
$ws.Range("C96").Value = -0.9225503716806988
$ws.Range("C97").Value = -0.8491698657853378
$ws.Range("D97").Value = 0.8526545954887239
$ws.Range("C98").Value = -1.191515643655161
$ws.Range("D98").Value = 1.016949629488724
$ws.Range("C99").Value = -1.430441087857995
$ws.Range("D99").Value = 0.8815376954887238
$ws.Range("B100").Value = -1.217909980957737
$ws.Range("C100").Value = 0.8721537754887239
$ws.Range("B101").Value = -1.26539928353432
$ws.Range("C101").Value = 0.5557457034887239
$ws.Range("C102").Value = 0.6390669014887239
$ws.Range("C103").Value = 0.5246197534887239
$ws.Range("D103").Value = 0.3082097950934801
$ws.Range("C104").Value = 0.2428672904887239
$ws.Range("D104").Value = 0.3048056840934801
$ws.Range("C105").Value = 0.5514316234887239
$ws.Range("D105").Value = 0.4343567680934801
$ws.Range("B106").Value = 0.3148484034887239
$ws.Range("C106").Value = 0.9884114590934802
$ws.Range("B107").Value = 0.3295457764887239
$ws.Range("C107").Value = 1.11054283609348
$ws.Range("C108").Value = 0.6209472400934801
$ws.Range("C109").Value = 0.9007734940934802
$ws.Range("D109").Value = 0.7305722247131936
$ws.Range("C110").Value = 0.5169714940934801
$ws.Range("D110").Value = 0.8421929667131937
$ws.Range("C111").Value = 0.2115338810934801
$ws.Range("D111").Value = 0.4788276257131937
$ws.Range("B112").Value = 0.0506311690934801
$ws.Range("C112").Value = -0.1957891702868064
$ws.Range("B113").Value = -0.02970344090651991
$ws.Range("C113").Value = -0.5800827092868064
$ws.Range("C114").Value = -0.2004644112868063
$ws.Range("C115").Value = -0.1799067152868063
$ws.Range("C116").Value = -0.4927747992868063
$ws.Range("C117").Value = -0.3419294162868063
$ws.Range("B119").Value = -0.3059158432868063
$ws.Range("B120").Value = -0.2670276532868063
$ws.Range("D122").Value = -0.4083682634916527
$ws.Range("D123").Value = -0.5334118554916527
$ws.Range("D124").Value = -0.4617206544916527
$ws.Range("D125").Value = -0.7651034874916527
$ws.Range("C126").Value = -0.2379622684916527
$ws.Range("C127").Value = -0.4135428994916527
$ws.Range("C128").Value = -0.4776197014916527
$ws.Range("C129").Value = -0.2036544774916527
$ws.Range("D129").Value = 0.5987332491758083
$ws.Range("C130").Value = -0.1998718194916527
$ws.Range("D130").Value = 0.6101592851758083
$ws.Range("C131").Value = -0.1834133014916527
$ws.Range("D131").Value = 0.5298148361758083
$ws.Range("C132").Value = 0.03791119950834732
$ws.Range("D132").Value = 0.5915373291758083
$ws.Range("B133").Value = -0.05564365149165268
$ws.Range("C133").Value = 0.4130432101758083
$ws.Range("B134").Value = 0.03353349250834725
$ws.Range("C134").Value = -0.3288647778241918
$ws.Range("C135").Value = 0.04982731217580827
$ws.Range("C136").Value = 0.3369026561758083
$ws.Range("D136").Value = 0.7254492243564907
$ws.Range("C137").Value = 0.3645469811758083
$ws.Range("D137").Value = 0.7215746373564907
$ws.Range("C138").Value = 0.1236970551758083
$ws.Range("D138").Value = 0.5311946523564907
$ws.Range("C139").Value = 0.2095627611758083
$ws.Range("D139").Value = 0.5539812373564907
$ws.Range("B140").Value = -0.08176241982419175
$ws.Range("C140").Value = 0.1753415943564907
$ws.Range("B141").Value = -0.1256759188241917
$ws.Range("C141").Value = 0.2651053283564908
$ws.Range("C142").Value = 0.08763596535649075
$ws.Range("C143").Value = 0.1003532183564907
$ws.Range("C144").Value = -0.02418658464350926
$ws.Range("C145").Value = 0.2001520573564908
